# ---------------------------------------------------------------------------
# cmms-pon.xlsx: "wo opening just working"
#
# 1. Insert a new "mode" sheet right after "action" (name/Forward/Reverse/Stay).
# 2. Insert a new "foreign_mode" column into the "action" sheet, classifying
#    each action as Forward / Reverse / Stay; existing "description" column
#    shifts one to the right.
# 3. Profile sheet: selection moves from B12 to D1.
# 4. action sheet: selection moves from C28 to B2 (stays the active tab).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Profile sheet: just a selection change -------------------------------
$profile = $wb.Worksheets.Item("Profile")
$profile.Range("D1").Select()

# --- action sheet: insert "foreign_mode" column ----------------------------
$action = $wb.Worksheets.Item("action")

# Shift the existing "description" column (B) one place to the right and
# make room for the new "foreign_mode" column.
$action.Columns("B:B").Insert()
$action.Columns("B:B").ColumnWidth = 16.67

# Header first, so "foreign_mode" lands in the shared-string table before
# the Forward/Reverse/Stay values that follow it.
$action.Range("B1").Value = "foreign_mode"

$action.Range("B2").Value = "Reverse"
$action.Range("B3").Value = "Forward"
$action.Range("B4").Value = "Reverse"
$action.Range("B5").Value = "Forward"
$action.Range("B6").Value = "Reverse"
$action.Range("B7").Value = "Forward"
$action.Range("B8").Value = "Reverse"
$action.Range("B9").Value = "Forward"
$action.Range("B10").Value = "Stay"
$action.Range("B11").Value = "Stay"
$action.Range("B12").Value = "Stay"
$action.Range("B13").Value = "Stay"
$action.Range("B14").Value = "Forward"
$action.Range("B15").Value = "Reverse"
$action.Range("B16").Value = "Reverse"
$action.Range("B17").Value = "Forward"
$action.Range("B18").Value = "Reverse"
$action.Range("B19").Value = "Stay"
$action.Range("B20").Value = "Reverse"

# --- New "mode" sheet, placed right after "action" -------------------------
$mode = $wb.Worksheets.Add($null, $action)
$mode.Name = "mode"
$mode.Tab.Color = 5287936
$mode.Columns("A:A").ColumnWidth = 16.67

$mode.Range("A1").Value = "name"
$mode.Range("A2").Value = "Forward"
$mode.Range("A3").Value = "Reverse"
$mode.Range("A4").Value = "Stay"
$mode.Range("A4").Select()

# --- Leave "action" as the active sheet/selection, matching the target -----
$action.Activate()
$action.Range("B2").Select()
